# "Fruta / hortaliza, semanal" -- weekly refresh of the Perejil
# (Vega Monumental Concepción) price series.
#
# The underlying data is organised as consecutive row-pairs (Primera /
# Segunda quality) each tagged with a single "Fecha" (column D). This
# week's refresh inserts two brand-new weekly observations (2021-08-27
# and 2021-08-25) into the series; every older observation's row-pair
# keeps its other columns but is shifted down to make room, and the two
# oldest observations that fall off the bottom of the existing range are
# appended as new rows 92-95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D ("Fecha") updates for the existing rows 50-91 ---------------
# (new value = the date that belonged to the row-pair one slot up, except
# at the two insertion points where a brand-new date is introduced)

$dateUpdates = @{
    50 = 44435; 51 = 44435
    52 = 44405; 53 = 44405
    54 = 44224; 55 = 44224
    56 = 44327; 57 = 44327
    58 = 44231; 59 = 44231
    60 = 44313; 61 = 44313
    62 = 44330; 63 = 44330
    64 = 44391; 65 = 44391
    66 = 44350; 67 = 44350
    68 = 44278; 69 = 44278
    70 = 44358; 71 = 44358
    72 = 44250; 73 = 44250
    74 = 44292; 75 = 44292
    76 = 44433; 77 = 44433
    78 = 44344; 79 = 44344
    80 = 44316; 81 = 44316
    82 = 44160; 83 = 44160
    84 = 44272; 85 = 44272
    86 = 44398; 87 = 44398
    88 = 44355; 89 = 44355
    90 = 44217; 91 = 44217
}

foreach ($row in $dateUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $dateUpdates[$row]
}

# --- Volumen (J) carried along with its row-pair for the 2021-03-23 entry -
# (this pair moved from rows 66/67 to rows 68/69)
$ws.Cells.Item(66, 10).Value = 200
$ws.Cells.Item(67, 10).Value = 100
$ws.Cells.Item(68, 10).Value = 300
$ws.Cells.Item(69, 10).Value = 150

# --- Origen (O) carried along with its row-pair ----------------------------
# Arica y Parinacota pair moved from rows 70/71 to rows 72/73
$ws.Cells.Item(70, 15).Value = "Región de Ñuble"
$ws.Cells.Item(71, 15).Value = "Región de Ñuble"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(73, 15).Value = "Región de Arica y Parinacota"

# --- Unidad de comercialización (N) carried along with its row-pair -------
# "$/docena de 1 kilo" pair moved from rows 74/75 to rows 78/79
$ws.Cells.Item(74, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(75, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(78, 14).Value = "$/docena de 1 kilo"
$ws.Cells.Item(79, 14).Value = "$/docena de 1 kilo"

# --- New rows 92-95: the two oldest observations pushed off the bottom ----

function Set-PerejilRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen
    )

    $ws.Cells.Item($Row, 1).Value = 11
    $ws.Cells.Item($Row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($Row, 3).Value = "Bíobío"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 8
    $ws.Cells.Item($Row, 6).Value = 100112044
    $ws.Cells.Item($Row, 7).Value = "Perejil"
    $ws.Cells.Item($Row, 8).Value = "Sin especificar"
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = 1
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# 2021-05-18 pair, pushed from rows 88/89 down to the new tail
Set-PerejilRow 92 44334 "Primera" 200 600 700 650 "$/atado 0,5 a 1 kilo" "Región de Ñuble"
Set-PerejilRow 93 44334 "Segunda" 100 500 500 500 "$/atado 0,5 a 1 kilo" "Región de Ñuble"

# 2020-12-03 pair, pushed from rows 90/91 down to the new tail
Set-PerejilRow 94 44168 "Primera" 200 600 700 650 "$/atado 0,5 a 1 kilo" "Región de Ñuble"
Set-PerejilRow 95 44168 "Segunda" 100 500 500 500 "$/atado 0,5 a 1 kilo" "Región de Ñuble"

"dimension now: " + $ws.UsedRange.Address()
